$d = $word.ActiveDocument

# --- Heading cleanup: merge "Struktura" + " " + "modela" into a single run ---
# (Text content is unchanged; this just collapses the pre-existing split runs
#  into one, matching the target's simplified markup.)
$r0 = $d.Content.Find.Execute("Struktura modela", $false, $false, $false, $false, $false, $true, 1, $false, "Struktura modela", 2)
Write-Output "heading merge: $r0"

# --- Expand the "Za čuvanje podataka..." (MySql/MongoDB) explanation paragraph ---
# Three insertions turn the short justification into the longer, more detailed
# explanation of why MySql and MongoDB were chosen.

# 1) Explain why MySql was chosen, right after "... kao što je mySql,".
$old1 = "kao što je mySql, dok je čuvanje"
$new1 = "kao što je mySql, zato što imamo malu količinu informacija o korisnicima koje čuvamo, koja je fiksna, i zato što je ova baza podataka open-source, I radi na različitim operativnim sistemima, dok je čuvanje"
$r1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "mysql reason insert: $r1"

# 2) Explain why MongoDB was chosen, right after "... kao što je mongoDB, koja".
$old2 = "kao što je mongoDB, koja na dosta bolji"
$new2 = "kao što je mongoDB, koja je takođe open-source, i koja na dosta bolji"
$r2 = $d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "mongodb reason insert: $r2"

# 3) Add "da ih često menjamo i da" before "njima pristupamo" after "veću količinu podataka,".
$old3 = "veću količinu podataka, da njima pristupamo"
$new3 = "veću količinu podataka, da ih često menjamo i da njima pristupamo"
$r3 = $d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "often-change insert: $r3"
